# Remake tables in Excel files
#
# raw: the "begin_date" column (D) used to be a literal "null" placeholder
#      for every row; it now mirrors the real sign_date (column B) value
#      (same text, same per-row number format).
# create-contracts: the begin_date column's SQL type suffix changes from
#      "," (nullable) to "NOT NULL,".
# insert-contracts: the generated INSERT statements stop quoting the
#      end_date value (raw!C, still the literal "null") so it is emitted
#      as a real SQL NULL instead of the string 'null'.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "raw": fill in D2:D5 with the same text + per-row format as B2:B5
# ---------------------------------------------------------------------
$raw = $wb.Worksheets.Item("raw")

$raw.Range("D2").Value = "'2016-09-23"
$raw.Range("B2").Copy() | Out-Null
$raw.Range("D2").PasteSpecial(-4122) | Out-Null

$raw.Range("D3").Value = "'2017-12-29"
$raw.Range("B3").Copy() | Out-Null
$raw.Range("D3").PasteSpecial(-4122) | Out-Null

$raw.Range("D4").Value = "'2017-11-23"
$raw.Range("B4").Copy() | Out-Null
$raw.Range("D4").PasteSpecial(-4122) | Out-Null

$raw.Range("D5").Value = "'2016-12-19"
$raw.Range("B5").Copy() | Out-Null
$raw.Range("D5").PasteSpecial(-4122) | Out-Null

$raw.Activate() | Out-Null
$raw.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "create-contracts": begin_date becomes NOT NULL
# ---------------------------------------------------------------------
$create = $wb.Worksheets.Item("create-contracts")
$create.Range("C5").Value = "NOT NULL,"

$create.Activate() | Out-Null
$create.Range("C5").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "insert-contracts": stop quoting the (still-null) end_date field
# ---------------------------------------------------------------------
$ins = $wb.Worksheets.Item("insert-contracts")

$ins.Range("A2").Formula = "=CONCATENATE(""INSERT INTO "",'create-contracts'!`$B`$1,"" VALUES ("",""'"",raw!A2,""', "",""'"",raw!B2,""', "",raw!C2,"", "",""'"",raw!D2,""', "",""'"",raw!E2,""', "",""'"",raw!F2,""', "",""'"",raw!G2,""');"")"
$ins.Range("A3").Formula = "=CONCATENATE(""INSERT INTO "",'create-contracts'!`$B`$1,"" VALUES ("",""'"",raw!A3,""', "",""'"",raw!B3,""', "",raw!C3,"", "",""'"",raw!D3,""', "",""'"",raw!E3,""', "",""'"",raw!F3,""', "",""'"",raw!G3,""');"")"
$ins.Range("A4").Formula = "=CONCATENATE(""INSERT INTO "",'create-contracts'!`$B`$1,"" VALUES ("",""'"",raw!A4,""', "",""'"",raw!B4,""', "",raw!C4,"", "",""'"",raw!D4,""', "",""'"",raw!E4,""', "",""'"",raw!F4,""', "",""'"",raw!G4,""');"")"
$ins.Range("A5").Formula = "=CONCATENATE(""INSERT INTO "",'create-contracts'!`$B`$1,"" VALUES ("",""'"",raw!A5,""', "",""'"",raw!B5,""', "",raw!C5,"", "",""'"",raw!D5,""', "",""'"",raw!E5,""', "",""'"",raw!F5,""', "",""'"",raw!G5,""');"")"

# widen the wrap-text cells in column A to also pin top-left alignment
$ins.Range("A2").HorizontalAlignment = -4131
$ins.Range("A2").VerticalAlignment = -4160
$ins.Range("A2").Copy() | Out-Null
$ins.Range("A3:A5").PasteSpecial(-4122) | Out-Null

$ins.Activate() | Out-Null
$ins.Range("A5").Select() | Out-Null
